$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column headers: E1 and F1
$ws.Range("E1").Value = "Prerequisites"
$ws.Range("F1").Value = "CoRequisites"

# Update selection to reflect the header row only
$ws.Range("A1:I1").Select()
